$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure columns so A:C lose their custom width (back to sheet
# default) while column D keeps an explicit (narrower) width. Deleting the
# original A:C shifts the old column D (with its width formatting) into A;
# inserting three fresh columns ahead of it then pushes it back out to D,
# leaving A:C with no width override and D still carrying a customizable
# width we set explicitly below.
$null = $ws.Columns("A:C").Delete()
$null = $ws.Columns.Item(1).Insert()
$null = $ws.Columns.Item(1).Insert()
$null = $ws.Columns.Item(1).Insert()

# --- Headers (row 1) ---
$ws.Cells.Item(1,1).Value = "X"
$ws.Cells.Item(1,2).Value = "Y"
$ws.Cells.Item(1,3).Value = "Radius"
$ws.Cells.Item(1,4).Value = "Comment"

# --- Clear out the old "Comment" column text for every data row, and any
# leftover data below the new data set ---
$ws.Range("A2:D200").ClearContents()

# --- New coordinate data (rows 2-17), no Comment entries this time ---
$ws.Cells.Item(2,1).Value = 156
$ws.Cells.Item(2,2).Value = -239
$ws.Cells.Item(2,3).Value = 8

$ws.Cells.Item(3,1).Value = 312
$ws.Cells.Item(3,2).Value = -239
$ws.Cells.Item(3,3).Value = 8

$ws.Cells.Item(4,1).Value = 238
$ws.Cells.Item(4,2).Value = -98
$ws.Cells.Item(4,3).Value = 8

$ws.Cells.Item(5,1).Value = 235
$ws.Cells.Item(5,2).Value = -83
$ws.Cells.Item(5,3).Value = 8

$ws.Cells.Item(6,1).Value = 320
$ws.Cells.Item(6,2).Value = -157
$ws.Cells.Item(6,3).Value = 8

$ws.Cells.Item(7,1).Value = 418
$ws.Cells.Item(7,2).Value = -287
$ws.Cells.Item(7,3).Value = 8

$ws.Cells.Item(8,1).Value = 299
$ws.Cells.Item(8,2).Value = -396
$ws.Cells.Item(8,3).Value = 8

$ws.Cells.Item(9,1).Value = 238
$ws.Cells.Item(9,2).Value = -514
$ws.Cells.Item(9,3).Value = 8

$ws.Cells.Item(10,1).Value = 238
$ws.Cells.Item(10,2).Value = -396
$ws.Cells.Item(10,3).Value = 1

$ws.Cells.Item(11,1).Value = 177
$ws.Cells.Item(11,2).Value = -396
$ws.Cells.Item(11,3).Value = 5

$ws.Cells.Item(12,1).Value = 299
$ws.Cells.Item(12,2).Value = -396
$ws.Cells.Item(12,3).Value = 5

$ws.Cells.Item(13,1).Value = 114
$ws.Cells.Item(13,2).Value = -344
$ws.Cells.Item(13,3).Value = 3

$ws.Cells.Item(14,1).Value = 363
$ws.Cells.Item(14,2).Value = -344
$ws.Cells.Item(14,3).Value = 3

$ws.Cells.Item(15,1).Value = 204
$ws.Cells.Item(15,2).Value = -277
$ws.Cells.Item(15,3).Value = 20

$ws.Cells.Item(16,1).Value = 271
$ws.Cells.Item(16,2).Value = -277
$ws.Cells.Item(16,3).Value = 20

$ws.Cells.Item(17,1).Value = 238
$ws.Cells.Item(17,2).Value = -298
$ws.Cells.Item(17,3).Value = 13

# --- Column D narrower width (~11.33 chars) ---
$ws.Columns.Item(4).ColumnWidth = 10.5

# --- Selection matches the saved file (D2 active) ---
$null = $ws.Range("D2").Select()
